# "Handles float input without breaking stuff"
#
# The marksheet previously showed an "Absent" student (no attempts, all
# zeros) with three side-by-side Student-Ans/Correct-Ans blocks (A:B,
# D:E, G:H). The student actually answered the quiz, so:
#   - the summary rows (10-12) get real Right/Wrong/Not-Attempted/Max
#     figures, a numeric (not text) marking penalty, and a real score
#   - the redundant 3rd block (G:H) and most of the 2nd block (D:E) are
#     removed - only two columns of answers (A:B) are kept, colour-coded
#     green/red per question depending on whether the student's answer
#     (column A) matches the correct answer (column B)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the unused 3rd (Student Ans / Correct Ans) block entirely, and the
# bulk of the 2nd block (only D16:E18 survive, handled individually below).
$ws.Range("G15:H21").Clear()
$ws.Range("D19:E40").Clear()

# --- Summary table (rows 10-12) ---
$ws.Range("A10").Value = "No."
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("A10").HorizontalAlignment = -4108

$ws.Range("B10").Value = 19
$ws.Range("B10").Style = "correctStyle"
$ws.Range("B10").HorizontalAlignment = -4108

$ws.Range("C10").Value = 2
$ws.Range("C10").Style = "incorrectStyle"
$ws.Range("C10").HorizontalAlignment = -4108

$ws.Range("D10").Value = 7
$ws.Range("D10").Style = "normalStyle"
$ws.Range("D10").HorizontalAlignment = -4108

$ws.Range("E10").Value = 28
$ws.Range("E10").Style = "normalStyle"
$ws.Range("E10").HorizontalAlignment = -4108

$ws.Range("A11").Value = "Marking"
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("A11").HorizontalAlignment = -4108

$ws.Range("B11").Value = 4
$ws.Range("B11").Style = "correctStyle"
$ws.Range("B11").HorizontalAlignment = -4108

# Previously stored as text "-1"; now a real number.
$ws.Range("C11").Value = -1
$ws.Range("C11").Style = "incorrectStyle"
$ws.Range("C11").HorizontalAlignment = -4108

$ws.Range("A12").Value = "Total"
$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("A12").HorizontalAlignment = -4108

$ws.Range("B12").Value = 76
$ws.Range("B12").Style = "correctStyle"
$ws.Range("B12").HorizontalAlignment = -4108

$ws.Range("C12").Value = -2
$ws.Range("C12").Style = "incorrectStyle"
$ws.Range("C12").HorizontalAlignment = -4108

$ws.Range("E12").Value = "74/112"
$ws.Range("E12").Style = "absoluteStyle"
$ws.Range("E12").HorizontalAlignment = -4108

# --- Remaining 2nd block entries (D16:E18 keep their original pairing) ---
$ws.Range("D16").Value = "Option A"
$ws.Range("D16").Style = "correctStyle"
$ws.Range("D16").HorizontalAlignment = -4108

$ws.Range("D17").Value = "Option A"
$ws.Range("D17").Style = "incorrectStyle"
$ws.Range("D17").HorizontalAlignment = -4108

# --- Student answers (column A), colour-coded against column B ---
$ws.Range("A18").Value = "Option B"
$ws.Range("A18").Style = "correctStyle"
$ws.Range("A18").HorizontalAlignment = -4108

$ws.Range("A19").Value = "Option C"
$ws.Range("A19").Style = "correctStyle"
$ws.Range("A19").HorizontalAlignment = -4108

$ws.Range("A20").Value = "Option B"
$ws.Range("A20").Style = "correctStyle"
$ws.Range("A20").HorizontalAlignment = -4108

$ws.Range("A21").Value = "Option C"
$ws.Range("A21").Style = "correctStyle"
$ws.Range("A21").HorizontalAlignment = -4108

$ws.Range("A22").Value = "Option D"
$ws.Range("A22").Style = "correctStyle"
$ws.Range("A22").HorizontalAlignment = -4108

$ws.Range("A23").Value = "Option D"
$ws.Range("A23").Style = "correctStyle"
$ws.Range("A23").HorizontalAlignment = -4108

$ws.Range("A24").Value = "Option A"
$ws.Range("A24").Style = "correctStyle"
$ws.Range("A24").HorizontalAlignment = -4108

$ws.Range("A25").Value = "Option A"
$ws.Range("A25").Style = "correctStyle"
$ws.Range("A25").HorizontalAlignment = -4108

$ws.Range("A26").Value = "Option C"
$ws.Range("A26").Style = "correctStyle"
$ws.Range("A26").HorizontalAlignment = -4108

$ws.Range("A27").Value = "Option A"
$ws.Range("A27").Style = "correctStyle"
$ws.Range("A27").HorizontalAlignment = -4108

$ws.Range("A28").Value = "Option B"
$ws.Range("A28").Style = "incorrectStyle"
$ws.Range("A28").HorizontalAlignment = -4108

# A29, A30, A34, A35 have no student answer (left blank / not attempted).

$ws.Range("A31").Value = "Option D"
$ws.Range("A31").Style = "correctStyle"
$ws.Range("A31").HorizontalAlignment = -4108

$ws.Range("A32").Value = "Option C"
$ws.Range("A32").Style = "correctStyle"
$ws.Range("A32").HorizontalAlignment = -4108

$ws.Range("A33").Value = "Option D"
$ws.Range("A33").Style = "correctStyle"
$ws.Range("A33").HorizontalAlignment = -4108

$ws.Range("A36").Value = "Option A"
$ws.Range("A36").Style = "correctStyle"
$ws.Range("A36").HorizontalAlignment = -4108

$ws.Range("A37").Value = "Option A"
$ws.Range("A37").Style = "correctStyle"
$ws.Range("A37").HorizontalAlignment = -4108

$ws.Range("A38").Value = "Option A"
$ws.Range("A38").Style = "correctStyle"
$ws.Range("A38").HorizontalAlignment = -4108

$ws.Range("A39").Value = "Option D"
$ws.Range("A39").Style = "correctStyle"
$ws.Range("A39").HorizontalAlignment = -4108

$ws.Range("A40").Value = "Option D"
$ws.Range("A40").Style = "correctStyle"
$ws.Range("A40").HorizontalAlignment = -4108
